# Apply the "added casefile trials tc, regular e2e tc" edit:
#   - the old "Ethnicity" header in column I is removed
#   - row 2 is populated with the result data for the casefile trial test case

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused "Ethnicity" header cell (column I)
$ws.Range("I1").ClearContents()

# Fill in the data row with the new trial/casefile values
$ws.Range("A2").Value2 = "CTDC-46730"
$ws.Range("B2").Value2 = "NCI-MATCH"
$ws.Range("C2").Value2 = "Q"
$ws.Range("D2").Value2 = "Ado-trastuzumab Emtansine"
$ws.Range("E2").Value2 = "Adenocarcinoma of the cervix"
$ws.Range("F2").Value2 = "FEMALE"
$ws.Range("G2").Value2 = "UNKNOWN"
$ws.Range("H2").Value2 = "UNKNOWN"
